# Generate Report for Handoff
# Replaces the two handed-back e2e source files with a new pair of files
# that are now "Ready for handoff", refreshes their generated artifact
# names/timestamps, and clears the stale "Latest Target File" / "Latest
# Handback File" columns (no handback has happened yet for the new files).

$wb = $excel.ActiveWorkbook

$oldFile1 = "764aa2b8-bf11-4cf0-8544-6c5e104ae578"
$oldFile2 = "c7aebfb3-e028-4625-a464-c4a202a9e3ed"
$newFile1 = "d4bc9dbf-0870-47a8-963f-5aede2dd074f"
$newFile2 = "ffffa1516d6c-1b7a-44e3-a722-b57c9951bffb"

$oldHash1 = "936d422752ada39154efa47a1f3bdfe180c42165"
$oldHash2 = "b4b684eb810b3565fcd59ccf816335e12673bde2"
$newHash  = "3d9250b61cab37357bbae60693ea6debfc0d761b"

$statusText = "Ready for handoff"

$mainRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41bd1bb8bf6c9694c7d6c7d0e5a5ec109b512580/e2e/"
$zhcnRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/80df63adbce4e30055cdc85cf419f970de1b4450/e2e/"
$dedeRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/057f9f47014f4a5358a7be4b3d1f3573cad21cdd/e2e/"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "$newFile1.md"
$wsOverview.Range("B2").Value = "e2e\$newFile1.md"
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = "2016-08-22 17:05:25"

$wsOverview.Range("A3").Value = "$newFile2.md"
$wsOverview.Range("B3").Value = "e2e\$newFile2.md"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = "2016-08-22 17:05:25"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($mainRepoBase + "$newFile1.md"), "", "", "e2e\$newFile1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($mainRepoBase + "$newFile2.md"), "", "", "e2e\$newFile2.md") | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()

$wsZhCn.Range("A2").Value = "$newFile1.md"
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("G2").Value = "$newFile1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 17:05:20"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Range("A3").Value = "$newFile2.md"
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "$newFile1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-22 17:05:20"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"

$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I3").Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), ($mainRepoBase + "$newFile1.md"), "", "", "$newFile1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($mainRepoBase + "$newFile2.md"), "", "", "$newFile2.md") | Out-Null

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()

$wsDeDe.Range("A2").Value = "$newFile1.md"
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("G2").Value = "$newFile1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-22 17:05:25"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("A3").Value = "$newFile2.md"
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "$newFile1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-22 17:05:25"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I3").Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), ($mainRepoBase + "$newFile1.md"), "", "", "$newFile1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($mainRepoBase + "$newFile2.md"), "", "", "$newFile2.md") | Out-Null

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
